$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Nov 12 18:04:45 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:04:59 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:05:12 EST 2024"
